$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Font color: change theme color to explicit black RGB for the whole used range
$ws.Range("A1:G15").Font.Color = 0

# 2. Row heights: rows 1-14 -> 19.5, row 15 -> 20.25
for ($r = 1; $r -le 14; $r++) {
    $ws.Rows.Item($r).RowHeight = 19.5
}
$ws.Rows.Item(15).RowHeight = 20.25

# 3. Update G10 value from 2 to 1
$ws.Range("G10").Value = 1
